$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.031151043873883
$ws.Cells.Item(2, 4).Value = 1.033794005858964
$ws.Cells.Item(2, 5).Value = 1.04004754849501
$ws.Cells.Item(2, 6).Value = 1.049295746533633
$ws.Cells.Item(2, 9).Value = 1.034252190235561
$ws.Cells.Item(2, 10).Value = 1.036288353266912
$ws.Cells.Item(2, 11).Value = 1.036595331243573
$ws.Cells.Item(2, 12).Value = 1.042831005245076
$ws.Cells.Item(2, 13).Value = 1.052053201275855

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.032114475006056
$ws.Cells.Item(3, 4).Value = 1.034486947339885
$ws.Cells.Item(3, 5).Value = 1.040911172603168
$ws.Cells.Item(3, 6).Value = 1.050263392983344
$ws.Cells.Item(3, 9).Value = 1.034434255854914
$ws.Cells.Item(3, 10).Value = 1.036893467496079
$ws.Cells.Item(3, 11).Value = 1.037097636706108
$ws.Cells.Item(3, 12).Value = 1.043504832478979
$ws.Cells.Item(3, 13).Value = 1.052832663760756

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.03273813800451
$ws.Cells.Item(4, 4).Value = 1.0349351517338
$ws.Cells.Item(4, 5).Value = 1.041470595795955
$ws.Cells.Item(4, 6).Value = 1.050890145098847
$ws.Cells.Item(4, 9).Value = 1.034550333668626
$ws.Cells.Item(4, 10).Value = 1.037284676409204
$ws.Cells.Item(4, 11).Value = 1.037421812309858
$ws.Cells.Item(4, 12).Value = 1.043940804187239
$ws.Cells.Item(4, 13).Value = 1.05333702831738

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.033000387040909
$ws.Cells.Item(5, 4).Value = 1.035123533593508
$ws.Cells.Item(5, 5).Value = 1.041705919885412
$ws.Cells.Item(5, 6).Value = 1.051153778930654
$ws.Cells.Item(5, 9).Value = 1.03459871775174
$ws.Cells.Item(5, 10).Value = 1.037449058230777
$ws.Cells.Item(5, 11).Value = 1.037557891334787
$ws.Cells.Item(5, 12).Value = 1.04412407623866
$ws.Cells.Item(5, 13).Value = 1.053549062128902

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.03304442336831
$ws.Cells.Item(6, 4).Value = 1.035155161159996
$ws.Cells.Item(6, 5).Value = 1.041745440161006
$ws.Cells.Item(6, 6).Value = 1.051198052809794
$ws.Cells.Item(6, 9).Value = 1.03460681729514
$ws.Cells.Item(6, 10).Value = 1.037476653803436
$ws.Cells.Item(6, 11).Value = 1.037580727586689
$ws.Cells.Item(6, 12).Value = 1.044154847780878
$ws.Cells.Item(6, 13).Value = 1.053584663445432

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.032741641949849
$ws.Cells.Item(7, 4).Value = 1.034937669073152
$ws.Cells.Item(7, 5).Value = 1.041473739649249
$ws.Cells.Item(7, 6).Value = 1.05089366720999
$ws.Cells.Item(7, 9).Value = 1.034550981810741
$ws.Cells.Item(7, 10).Value = 1.037286873212081
$ws.Cells.Item(7, 11).Value = 1.03742363140886
$ws.Cells.Item(7, 12).Value = 1.043943253119677
$ws.Cells.Item(7, 13).Value = 1.053339861529049

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.031476585899872
$ws.Cells.Item(8, 4).Value = 1.034028224116637
$ws.Cells.Item(8, 5).Value = 1.040339289207547
$ws.Cells.Item(8, 6).Value = 1.049622638038526
$ws.Cells.Item(8, 9).Value = 1.034314078140005
$ws.Cells.Item(8, 10).Value = 1.036492924407443
$ws.Cells.Item(8, 11).Value = 1.036765263046733
$ws.Cells.Item(8, 12).Value = 1.043058735917389
$ws.Cells.Item(8, 13).Value = 1.052316623569898

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.029249411610389
$ws.Cells.Item(9, 4).Value = 1.032424384396941
$ws.Cells.Item(9, 5).Value = 1.03834489629044
$ws.Cells.Item(9, 6).Value = 1.047387725712584
$ws.Cells.Item(9, 9).Value = 1.033883396907239
$ws.Cells.Item(9, 10).Value = 1.035091317150397
$ws.Cells.Item(9, 11).Value = 1.035598670581238
$ws.Cells.Item(9, 12).Value = 1.041499843641414
$ws.Cells.Item(9, 13).Value = 1.05051359734007

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.027766029879014
$ws.Cells.Item(10, 4).Value = 1.031354375255896
$ws.Cells.Item(10, 5).Value = 1.037018500751338
$ws.Cells.Item(10, 6).Value = 1.045901084762932
$ws.Cells.Item(10, 9).Value = 1.033587421932293
$ws.Cells.Item(10, 10).Value = 1.034155237735369
$ws.Cells.Item(10, 11).Value = 1.034816654629639
$ws.Cells.Item(10, 12).Value = 1.040460464285104
$ws.Cells.Item(10, 13).Value = 1.049311676951598

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.027124050529029
$ws.Cells.Item(11, 4).Value = 1.030890880147991
$ws.Cells.Item(11, 5).Value = 1.03644493032455
$ws.Cells.Item(11, 6).Value = 1.045258150174917
$ws.Cells.Item(11, 9).Value = 1.033457169325331
$ws.Cells.Item(11, 10).Value = 1.033749518414146
$ws.Cells.Item(11, 11).Value = 1.034477027976454
$ws.Cells.Item(11, 12).Value = 1.040010385232444
$ws.Cells.Item(11, 13).Value = 1.048791268166728

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.026885641751821
$ws.Cells.Item(12, 4).Value = 1.030718692546047
$ws.Cells.Item(12, 5).Value = 1.036231997184044
$ws.Cells.Item(12, 6).Value = 1.045019455560543
$ws.Cells.Item(12, 9).Value = 1.033408473775168
$ws.Cells.Item(12, 10).Value = 1.033598758291029
$ws.Cells.Item(12, 11).Value = 1.034350724893693
$ws.Cells.Item(12, 12).Value = 1.039843203580029
$ws.Cells.Item(12, 13).Value = 1.048597970775457

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.026936778928382
$ws.Cells.Item(13, 4).Value = 1.030755628476741
$ws.Cells.Item(13, 5).Value = 1.036277666777558
$ws.Cells.Item(13, 6).Value = 1.045070650921134
$ws.Cells.Item(13, 9).Value = 1.033418933332629
$ws.Cells.Item(13, 10).Value = 1.033631099456528
$ws.Cells.Item(13, 11).Value = 1.034377824135134
$ws.Cells.Item(13, 12).Value = 1.039879064699609
$ws.Cells.Item(13, 13).Value = 1.04863943345258

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.02710434254227
$ws.Cells.Item(14, 4).Value = 1.030876647567339
$ws.Cells.Item(14, 5).Value = 1.036427326822564
$ws.Cells.Item(14, 6).Value = 1.045238417150885
$ws.Cells.Item(14, 9).Value = 1.033453150536684
$ws.Cells.Item(14, 10).Value = 1.03373705771189
$ws.Cells.Item(14, 11).Value = 1.034466590791644
$ws.Cells.Item(14, 12).Value = 1.0399965659857
$ws.Cells.Item(14, 13).Value = 1.048775290024033

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.027207590784213
$ws.Cells.Item(15, 4).Value = 1.030951208164099
$ws.Cells.Item(15, 5).Value = 1.036519552780825
$ws.Cells.Item(15, 6).Value = 1.045341799385721
$ws.Cells.Item(15, 9).Value = 1.03347419130394
$ws.Cells.Item(15, 10).Value = 1.033802334441454
$ws.Cells.Item(15, 11).Value = 1.034521262926581
$ws.Cells.Item(15, 12).Value = 1.040068962134274
$ws.Cells.Item(15, 13).Value = 1.04885899651364

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.027808643041824
$ws.Cells.Item(16, 4).Value = 1.031385132340677
$ws.Cells.Item(16, 5).Value = 1.03705658300056
$ws.Cells.Item(16, 6).Value = 1.045943770969697
$ws.Cells.Item(16, 9).Value = 1.033596022311906
$ws.Cells.Item(16, 10).Value = 1.034182155826411
$ws.Cells.Item(16, 11).Value = 1.034839173359824
$ws.Cells.Item(16, 12).Value = 1.040490334184727
$ws.Cells.Item(16, 13).Value = 1.049346215499911

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.028185757317251
$ws.Cells.Item(17, 4).Value = 1.03165727580466
$ws.Cells.Item(17, 5).Value = 1.037393654175992
$ws.Cells.Item(17, 6).Value = 1.046321584336492
$ws.Cells.Item(17, 9).Value = 1.033671883615792
$ws.Cells.Item(17, 10).Value = 1.034420303698185
$ws.Cells.Item(17, 11).Value = 1.035038320789926
$ws.Cells.Item(17, 12).Value = 1.040754644765972
$ws.Cells.Item(17, 13).Value = 1.049651844012913

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.028405753946774
$ws.Cells.Item(18, 4).Value = 1.031815995600804
$ws.Cells.Item(18, 5).Value = 1.037590336075822
$ws.Cells.Item(18, 6).Value = 1.04654203270274
$ws.Cells.Item(18, 9).Value = 1.033715930200028
$ws.Cells.Item(18, 10).Value = 1.034559173566425
$ws.Cells.Item(18, 11).Value = 1.035154382742135
$ws.Cells.Item(18, 12).Value = 1.04090881057683
$ws.Cells.Item(18, 13).Value = 1.049830114846862

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.028480772563165
$ws.Cells.Item(19, 4).Value = 1.031870112047285
$ws.Cells.Item(19, 5).Value = 1.037657412035783
$ws.Cells.Item(19, 6).Value = 1.046617212793225
$ws.Cells.Item(19, 9).Value = 1.033730914661037
$ws.Cells.Item(19, 10).Value = 1.034606518174163
$ws.Cells.Item(19, 11).Value = 1.035193940324487
$ws.Cells.Item(19, 12).Value = 1.040961376754447
$ws.Cells.Item(19, 13).Value = 1.049890901054757

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.028145293183642
$ws.Cells.Item(20, 4).Value = 1.031628079123467
$ws.Cells.Item(20, 5).Value = 1.037357481967937
$ws.Cells.Item(20, 6).Value = 1.046281040645757
$ws.Cells.Item(20, 9).Value = 1.03366376530853
$ws.Cells.Item(20, 10).Value = 1.0343947565848
$ws.Cells.Item(20, 11).Value = 1.035016964220691
$ws.Cells.Item(20, 12).Value = 1.040726286950606
$ws.Cells.Item(20, 13).Value = 1.049619052659423

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.027054997823055
$ws.Cells.Item(21, 4).Value = 1.030841011138141
$ws.Cells.Item(21, 5).Value = 1.036383252445704
$ws.Cells.Item(21, 6).Value = 1.045189010858397
$ws.Cells.Item(21, 9).Value = 1.033443083084273
$ws.Cells.Item(21, 10).Value = 1.033705857238571
$ws.Cells.Item(21, 11).Value = 1.034440455374392
$ws.Cells.Item(21, 12).Value = 1.039961964840611
$ws.Cells.Item(21, 13).Value = 1.048735283487242

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.026369780886735
$ws.Cells.Item(22, 4).Value = 1.030346006925102
$ws.Cells.Item(22, 5).Value = 1.035771390004776
$ws.Cells.Item(22, 6).Value = 1.044503102750176
$ws.Cells.Item(22, 9).Value = 1.033302515385199
$ws.Cells.Item(22, 10).Value = 1.033272384267554
$ws.Cells.Item(22, 11).Value = 1.034077110375351
$ws.Cells.Item(22, 12).Value = 1.039481392602298
$ws.Cells.Item(22, 13).Value = 1.048179655637715

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.026732999167153
$ws.Cells.Item(23, 4).Value = 1.030608431186855
$ws.Cells.Item(23, 5).Value = 1.036095685565885
$ws.Cells.Item(23, 6).Value = 1.044866649426369
$ws.Cells.Item(23, 9).Value = 1.033377204923862
$ws.Cells.Item(23, 10).Value = 1.033502207940211
$ws.Cells.Item(23, 11).Value = 1.034269808693642
$ws.Cells.Item(23, 12).Value = 1.039736153934704
$ws.Cells.Item(23, 13).Value = 1.048474201023128

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.028163577084192
$ws.Cells.Item(24, 4).Value = 1.031641271899213
$ws.Cells.Item(24, 5).Value = 1.037373826403223
$ws.Cells.Item(24, 6).Value = 1.046299360358974
$ws.Cells.Item(24, 9).Value = 1.033667434246409
$ws.Cells.Item(24, 10).Value = 1.034406300342378
$ws.Cells.Item(24, 11).Value = 1.035026614635975
$ws.Cells.Item(24, 12).Value = 1.040739100633202
$ws.Cells.Item(24, 13).Value = 1.049633869651696

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.029824945931131
$ws.Cells.Item(25, 4).Value = 1.032839159114891
$ws.Cells.Item(25, 5).Value = 1.038859935865062
$ws.Cells.Item(25, 6).Value = 1.047964927556434
$ws.Cells.Item(25, 9).Value = 1.033996301757562
$ws.Cells.Item(25, 10).Value = 1.035453965176408
$ws.Cells.Item(25, 11).Value = 1.035901022645645
$ws.Cells.Item(25, 12).Value = 1.041902879590032
$ws.Cells.Item(25, 13).Value = 1.050979710372256
